$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C holds the "Förändrad" (Changed) date. Every populated row in that
# column needs to move forward by one day (45188 -> 45189).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) {
    $lastRow = 351
}

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $cur = $cell.Value()
    if ($cur -ne $null) {
        $cell.Value = $cur.AddDays(1)
    }
}
